$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 (ano 2025) with refreshed metrics
$ws.Range("C6").Value = 426
$ws.Range("E6").Value = 117
$ws.Range("G6").Value = 27.46478873239437
$ws.Range("H6").Value = 72.53521126760563
